$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Policies")

# Insert a new column before column I (9th column: "isParent") to make room for
# the new "EarlyDuration" column. This shifts isParent -> J and RunMode -> K.
$ws.Columns.Item(9).Insert()

# Restore the column width of the newly inserted column to match its neighbors
# (this empirically maps to ~11.52 character-width units once saved).
$ws.Columns.Item(9).ColumnWidth = 10.6867

# Header for the new column
$ws.Range("I1").Value = "EarlyDuration"

# Data rows default to "no" for the new EarlyDuration column
$ws.Range("I2:I34").Value = "no"

# Refresh the AutoFilter so it covers the new column range A1:K34 (was A1:J34)
$ws.AutoFilterMode = $false
$ws.Range("A1:K34").AutoFilter() | Out-Null

# Update the hidden filter-database defined names that track the filter range,
# which otherwise keep pointing at the old J-column boundary.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Policies!_FilterDatabase") {
        $n.RefersTo = "=Policies!`$A`$1:`$K`$34"
    } elseif ($n.Name -eq "Policies!_FilterDatabase_0") {
        $n.RefersTo = "=Policies!`$A`$1:`$K`$18"
    } elseif ($n.Name -eq "Policies!_FilterDatabase_0_0") {
        $n.RefersTo = "=Policies!`$A`$1:`$K`$18"
    }
}

# Reflect the final cursor/selection position on the Policies sheet
$ws.Range("A18").Select() | Out-Null
